# Fruta / hortaliza, semanal
# Insert 6 new weekly price rows (newest week, date 44931) at the top of the
# "Femacal de La Calera - Cereza" data block (rows 678..683), pushing the
# existing rows down by 6 (old row 678 becomes row 684, ..., old row 713
# becomes row 719). Row 1..677 are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 678:713 down by inserting 6 blank rows above the current row 678.
$ws.Rows("678:683").Insert()

# Common fields shared by all new rows in this block.
$mercadoId = 3
$mercado   = "Femacal de La Calera"
$region    = "Coquimbo"
$fecha     = 44931
$codreg    = 5
$tipo      = "Fruta"
$prodId    = 100103
$producto  = "Frutos de hueso (carozo)"
$catId     = 100103001
$categoria = "Cereza"
$unidad    = "$/bandeja 10 kilos"
$origen    = "Provincia de Curicó"
$kgUnidad  = 10

# Variety, quality, volume, min, max, avg, price/kg for the 6 new rows.
$newRows = @(
    @(678, "Bing",        "Especial", 56, 7000, 7000, 7000, 700),
    @(679, "Bing",        "Primera",  60, 6000, 6000, 6000, 600),
    @(680, "Bing",        "Segunda",  60, 5000, 5000, 5000, 500),
    @(681, "Sweet Heart",  "Especial", 54, 7000, 7000, 7000, 700),
    @(682, "Sweet Heart",  "Primera",  57, 6000, 6000, 6000, 600),
    @(683, "Sweet Heart",  "Segunda",  50, 5000, 5000, 5000, 500)
)

foreach ($row in $newRows) {
    $r         = $row[0]
    $variedad  = $row[1]
    $calidad   = $row[2]
    $volumen   = $row[3]
    $precioMin = $row[4]
    $precioMax = $row[5]
    $precioProm = $row[6]
    $precioKg  = $row[7]

    $ws.Cells.Item($r, 1).Value  = $mercadoId
    $ws.Cells.Item($r, 2).Value  = $mercado
    $ws.Cells.Item($r, 3).Value  = $region
    $ws.Cells.Item($r, 4).Value  = $fecha
    $ws.Cells.Item($r, 5).Value  = $codreg
    $ws.Cells.Item($r, 6).Value  = $tipo
    $ws.Cells.Item($r, 7).Value  = $prodId
    $ws.Cells.Item($r, 8).Value  = $producto
    $ws.Cells.Item($r, 9).Value  = $catId
    $ws.Cells.Item($r, 10).Value = $categoria
    $ws.Cells.Item($r, 11).Value = $variedad
    $ws.Cells.Item($r, 12).Value = $calidad
    $ws.Cells.Item($r, 13).Value = $volumen
    $ws.Cells.Item($r, 14).Value = $precioMin
    $ws.Cells.Item($r, 15).Value = $precioMax
    $ws.Cells.Item($r, 16).Value = $precioProm
    $ws.Cells.Item($r, 17).Value = $unidad
    $ws.Cells.Item($r, 18).Value = $origen
    $ws.Cells.Item($r, 19).Value = $precioKg
    $ws.Cells.Item($r, 20).Value = $kgUnidad
}
